# Saldo.xlsx ("Export" sheet) update:
#   - add account 004574428 / GUILHERME / 300000 as the first data row
#   - add account 001761119 / BLUEMETRIX / 570.71 just above account 004342617 / JURACI
#   - remove account 004237325 / RICARDO (balance 332.66)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert GUILHERME as the new first data row (row 2, right under the header) ---
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = "'004574428"
$ws.Cells.Item(2, 2).Value = "GUILHERME"
$ws.Cells.Item(2, 3).Value = 300000

# --- 2) Insert BLUEMETRIX immediately above the JURACI row (now shifted to row 60) ---
$ws.Rows.Item(60).Insert()
$ws.Cells.Item(60, 1).Value = "'001761119"
$ws.Cells.Item(60, 2).Value = "BLUEMETRIX"
$ws.Cells.Item(60, 3).Value = 570.71

# --- 3) Delete the RICARDO row (now shifted to row 102) ---
$ws.Rows.Item(102).Delete()
